$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93. This shifts the existing rows 93-150
# down to 94-151, preserving all of their original data untouched.
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new data record.
$ws.Cells.Item(93, 1).Value = 4
$ws.Cells.Item(93, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(93, 3).Value = "Los Lagos"
$ws.Cells.Item(93, 4).Value = 44518
$ws.Cells.Item(93, 5).Value = 10
$ws.Cells.Item(93, 6).Value = 100112032
$ws.Cells.Item(93, 7).Value = "Zapallo italiano"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 120
$ws.Cells.Item(93, 11).Value = 12000
$ws.Cells.Item(93, 12).Value = 12000
$ws.Cells.Item(93, 13).Value = 12000
$ws.Cells.Item(93, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(93, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(93, 16).Value = 240
$ws.Cells.Item(93, 17).Value = 50
$ws.Cells.Item(93, 18).Value = "Hortaliza"
